$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2201942
$ws.Range("I69").Value = 4500
$ws.Range("J69").Value = 2401709.2
$ws.Range("K69").Value = 13500
$ws.Range("L69").Value = 7205127.600000001
$ws.Range("M69").Value = -12626
$ws.Range("N69").Value = -7206875.600000001

$ws.Range("H72").Value = 2201942
$ws.Range("I72").Value = 4500
$ws.Range("J72").Value = 2401709.2
$ws.Range("K72").Value = 40500
$ws.Range("L72").Value = 21615382.8
$ws.Range("M72").Value = -36132
$ws.Range("N72").Value = -21624118.8

$ws.Range("H129").Value = 1092.9066
$ws.Range("J129").Value = 1113.4445
$ws.Range("L129").Value = 3340.3335
$ws.Range("N129").Value = -13340.3335


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 308929.8
$ws.Range("I61").Value = 7564.737
$ws.Range("J61").Value = 717925.3
$ws.Range("K61").Value = 7564.737
$ws.Range("L61").Value = 717925.3
$ws.Range("M61").Value = -7352.737
$ws.Range("N61").Value = -718349.3

$ws.Range("H136").Value = 308929.8
$ws.Range("I136").Value = 7564.737
$ws.Range("J136").Value = 717925.3
$ws.Range("K136").Value = 22694.211
$ws.Range("L136").Value = 2153775.9
$ws.Range("M136").Value = -20144.211
$ws.Range("N136").Value = -2158875.9


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15702
$ws.Range("I20").Value = 1763.6
$ws.Range("J20").Value = 33125
$ws.Range("K20").Value = 1763.6
$ws.Range("L20").Value = 33125
$ws.Range("M20").Value = -1516.6
$ws.Range("N20").Value = -33619

$ws.Range("H75").Value = 17594.715
$ws.Range("I75").Value = 9106.75
$ws.Range("J75").Value = 28912
$ws.Range("K75").Value = 9106.75
$ws.Range("L75").Value = 28912
$ws.Range("M75").Value = -8170.75
$ws.Range("N75").Value = -30784

$ws.Range("H78").Value = 17594.715
$ws.Range("I78").Value = 9106.75
$ws.Range("J78").Value = 28912
$ws.Range("K78").Value = 27320.25
$ws.Range("L78").Value = 86736
$ws.Range("M78").Value = -22640.25
$ws.Range("N78").Value = -96096

$ws.Range("H99").Value = 41668108
$ws.Range("I99").Value = 55556892
$ws.Range("K99").Value = 55556892
$ws.Range("M99").Value = -55555394

$ws.Range("H134").Value = 24321.123
$ws.Range("I134").Value = 4465.8857
$ws.Range("J134").Value = 73959.21000000001
$ws.Range("K134").Value = 13397.6571
$ws.Range("L134").Value = 221877.63
$ws.Range("M134").Value = -10862.6571
$ws.Range("N134").Value = -226947.63


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 557724.9399999999
$ws.Range("I58").Value = 2080.6155
$ws.Range("J58").Value = 2002400.2
$ws.Range("K58").Value = 2080.6155
$ws.Range("L58").Value = 2002400.2
$ws.Range("M58").Value = -1877.6155
$ws.Range("N58").Value = -2002806.2

$ws.Range("H94").Value = 7534.467
$ws.Range("I94").Value = 6902.5713
$ws.Range("J94").Value = 8087.375
$ws.Range("K94").Value = 6902.5713
$ws.Range("L94").Value = 8087.375
$ws.Range("M94").Value = -6451.5713
$ws.Range("N94").Value = -8989.375

$ws.Range("H105").Value = 3467.6191
$ws.Range("I105").Value = 3472.1052
$ws.Range("J105").Value = 3425
$ws.Range("K105").Value = 3472.1052
$ws.Range("L105").Value = 3425
$ws.Range("M105").Value = -1725.1052
$ws.Range("N105").Value = -6919

$ws.Range("H136").Value = 557724.9399999999
$ws.Range("I136").Value = 2080.6155
$ws.Range("J136").Value = 2002400.2
$ws.Range("K136").Value = 6241.8465
$ws.Range("L136").Value = 6007200.6
$ws.Range("M136").Value = -3691.8465
$ws.Range("N136").Value = -6012300.6


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5455172
$ws.Range("I113").Value = 8333983
$ws.Range("J113").Value = 2000599.4
$ws.Range("K113").Value = 25001949
$ws.Range("L113").Value = 6001798.199999999
$ws.Range("M113").Value = -24999779
$ws.Range("N113").Value = -6006138.199999999

$ws.Range("H129").Value = 22223312
$ws.Range("I129").Value = 111111480
$ws.Range("J129").Value = 1269.4166
$ws.Range("K129").Value = 333334440
$ws.Range("L129").Value = 3808.2498
$ws.Range("M129").Value = -333329440
$ws.Range("N129").Value = -13808.2498

$ws.Range("H131").Value = 2128579.2
$ws.Range("I131").Value = 7143305.5
$ws.Range("J131").Value = 1119.5151
$ws.Range("K131").Value = 21429916.5
$ws.Range("L131").Value = 3358.5453
$ws.Range("M131").Value = -21424876.5
$ws.Range("N131").Value = -13438.5453


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5579.75
$ws.Range("I70").Value = 5644.3887
$ws.Range("J70").Value = 4998
$ws.Range("K70").Value = 5644.3887
$ws.Range("L70").Value = 4998
$ws.Range("M70").Value = -5374.3887
$ws.Range("N70").Value = -5538

$ws.Range("H73").Value = 5579.75
$ws.Range("I73").Value = 5644.3887
$ws.Range("J73").Value = 4998
$ws.Range("K73").Value = 5644.3887
$ws.Range("L73").Value = 4998
$ws.Range("M73").Value = -4708.3887
$ws.Range("N73").Value = -6870

$ws.Range("H80").Value = 6153.4614
$ws.Range("I80").Value = 9991.154
$ws.Range("J80").Value = 2315.7693
$ws.Range("K80").Value = 9991.154
$ws.Range("L80").Value = 2315.7693
$ws.Range("M80").Value = -8993.154
$ws.Range("N80").Value = -4311.7693

$ws.Range("H83").Value = 6153.4614
$ws.Range("I83").Value = 9991.154
$ws.Range("J83").Value = 2315.7693
$ws.Range("K83").Value = 49955.77
$ws.Range("L83").Value = 11578.8465
$ws.Range("M83").Value = -44963.77
$ws.Range("N83").Value = -21562.8465

$ws.Range("H126").Value = 10069.542
$ws.Range("I126").Value = 12984.059
$ws.Range("J126").Value = 2991.4285
$ws.Range("K126").Value = 38952.177
$ws.Range("L126").Value = 8974.2855
$ws.Range("M126").Value = -36482.177
$ws.Range("N126").Value = -13914.2855

$ws.Range("H132").Value = 3930.7273
$ws.Range("I132").Value = 4448.7144
$ws.Range("J132").Value = 3024.25
$ws.Range("K132").Value = 13346.1432
$ws.Range("L132").Value = 9072.75
$ws.Range("M132").Value = -10816.1432
$ws.Range("N132").Value = -14132.75


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1135.7142
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1191.6666
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1191.6666
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -1781.6666

$ws.Range("H27").Value = 1135.7142
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1191.6666
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 1191.6666
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1405.6666

$ws.Range("H40").Value = 35716830
$ws.Range("I40").Value = 47621668
$ws.Range("J40").Value = 2314.1428
$ws.Range("K40").Value = 47621668
$ws.Range("L40").Value = 2314.1428
$ws.Range("M40").Value = -47621532
$ws.Range("N40").Value = -2586.1428

$ws.Range("H74").Value = 22500
$ws.Range("J74").Value = 22500
$ws.Range("L74").Value = 22500
$ws.Range("N74").Value = -24496

$ws.Range("H77").Value = 22500
$ws.Range("J77").Value = 22500
$ws.Range("L77").Value = 67500
$ws.Range("N77").Value = -77484

$ws.Range("H82").Value = 619854.6
$ws.Range("I82").Value = 1001968
$ws.Range("K82").Value = 1001968
$ws.Range("M82").Value = -1001607

$ws.Range("H85").Value = 619854.6
$ws.Range("I85").Value = 1001968
$ws.Range("K85").Value = 1001968
$ws.Range("M85").Value = -1000720

$ws.Range("H122").Value = 2332136
$ws.Range("I122").Value = 2982594.8
$ws.Range("J122").Value = 912953.6
$ws.Range("K122").Value = 8947784.399999999
$ws.Range("L122").Value = 2738860.8
$ws.Range("M122").Value = -8945334.399999999
$ws.Range("N122").Value = -2743760.8

$ws.Range("H132").Value = 7097368.5
$ws.Range("I132").Value = 9809577
$ws.Range("J132").Value = 3898.3076
$ws.Range("K132").Value = 29428731
$ws.Range("L132").Value = 11694.9228
$ws.Range("M132").Value = -29426201
$ws.Range("N132").Value = -16754.9228


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1414.434
$ws.Range("I132").Value = 943.56757
$ws.Range("J132").Value = 2503.3125
$ws.Range("K132").Value = 2830.70271
$ws.Range("L132").Value = 7509.9375
$ws.Range("M132").Value = -300.70271
$ws.Range("N132").Value = -12569.9375

